# SCD0018-015 - Penyelia, CRO CRM, dan SRM mengajukan data Non Sales.xlsx
# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from SCD0292 -> SCD0018
$ws.Name = "SCD0018"

# Update the TC_ID column (B) for the three data rows: DGS-307 -> SCD0018-015
$ws.Range("B2").Value = "SCD0018-015"
$ws.Range("B3").Value = "SCD0018-015"
$ws.Range("B4").Value = "SCD0018-015"

# Widen column B so the longer TC_ID text fits (was 9 chars wide, now ~12.57)
$ws.Columns("B").ColumnWidth = 11.7

# Move the active selection from O2 to B5
[void]$ws.Range("B5").Select()
